$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Repull data / push all data / mean calculation: update dSF (column F) values
$ws.Range("F6").Value = 0
$ws.Range("F9").Value = 2
$ws.Range("F12").Value = 0
$ws.Range("F15").Value = 2
